$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date that was bumped by one day
# (2023-10-06 -> 2023-10-07, i.e. serial 45205 -> 45206) for every
# data row from row 2 through row 61.
for ($row = 2; $row -le 61; $row++) {
    $ws.Cells.Item($row, 3).Value = 45206
}
